$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the title (2nd paragraph of the document).
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2. Swap the trailing "Prompt: ..." paragraph's text for the meta
#    description text (keeping its existing italic run formatting).
# ---------------------------------------------------------------------------
$oldPrompt = "Prompt: Create a cartoon-style feature image for the game " + [char]34 + "Big Win 777" + [char]34 + " that features a happy Maya warrior with glasses. Design specifications: - The image should have a bright and colorful background that represents the Las Vegas-style slot machine setting. - The Maya warrior should be the central focus of the image, with a big smile on their face and glasses to represent the modern touch to the game. - The warrior should be holding a diamond and a wheel of fortune to represent the game's Wild and Scatter symbols. - The image should have a playful and fun tone, appealing to players who enjoy traditional slot machines with a modern twist."
$newPrompt = "Explore the vintage atmosphere of Las Vegas with Big Win 777 slot game. Play for free with high payout percentages, traditional symbols, and bonus features."

$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Insert a new bold paragraph ("Play Big Win 777 Free - Retro Slot
#    Machine by Play N Go") right before that last paragraph.
#
#    InsertXML on a collapsed insertion point merges the *last* <w:p> of the
#    inserted fragment with whatever content already starts at that point,
#    so we feed it a duplicate copy of the (already updated) last
#    paragraph's text and then collapse the resulting duplicate text back
#    down to a single copy.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xmlFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Big Win 777 Free - Retro Slot Machine by Play N Go</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>' + $newPrompt + '</w:t></w:r></w:p>'
$insertionPoint.InsertXML($xmlFragment)

$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$duplicated = $newPrompt + $newPrompt
$finalPara.Range.Find.Execute($duplicated, $true, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 2) | Out-Null
